$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.030.43"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "2.733.24"
$ws.Range("E3").Value = "  -6.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "507.47"
$ws.Range("E5").Value = "  -3.72%  "
$ws.Range("D6").Value = "141.40"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").Value = "2.743.37"
$ws.Range("E9").Value = "  -5.90%  "
$ws.Range("D10").Value = "6.10"
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "0.350"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "3.211.05"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("D15").Value = "58.950.28"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "21.82"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "2.745.10"
$ws.Range("E18").Value = "  -5.21%  "
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").Value = "345.34"
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "63.21"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "7.53"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "19.20"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").Value = "149.62"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "5.44"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "0.962"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "36.17"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").Value = "2.184.25"
$ws.Range("E42").Value = "  -6.55%  "
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -6.37%  "
$ws.Range("D46").Value = "19.20"
$ws.Range("E46").Value = "  -7.61%  "
$ws.Range("D47").Value = "4.80"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "0.0887"
$ws.Range("E50").Value = "  -4.18%  "
$ws.Range("D51").Value = "18.20"
$ws.Range("E51").Value = "  -0.26%  "
